# Update countries & provincias Spain
# Applies the 18-Abril-2020 data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp in A1
#  - Refreshes case counters for Rumania, Uzbekistan and Etiopia (no reordering)
#  - Zambia's case count overtakes Cabo Verde/Polinesia Francesa/Uganda/Bahamas,
#    so those 5 rows are re-sorted (descending by Casos totales) -> rows 154..158
#  - Fiyi / Islas Virgenes de los Estados Unidos swap order (tied case count)
#  - Montserrat overtakes Islas Turcas y Caicos -> rows 195/196 swap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 18 de Abril de 2020 a las 12:23"

# --- Simple numeric refreshes (country stays in place) ----------------------
# Rumania (row 31)
$ws.Cells.Item(31, 1).Value = "Rumania"
$ws.Cells.Item(31, 2).Value = 8418
$ws.Cells.Item(31, 3).Value = 351
$ws.Cells.Item(31, 4).Value = 1730
$ws.Cells.Item(31, 5).Value = 6271
$ws.Cells.Item(31, 6).Value = 248
$ws.Cells.Item(31, 7).Value = 6
$ws.Cells.Item(31, 8).Value = 417

# Uzbekistan (row 70)
$ws.Cells.Item(70, 1).Value = "Uzbekistan"
$ws.Cells.Item(70, 2).Value = 1450
$ws.Cells.Item(70, 3).Value = 45
$ws.Cells.Item(70, 4).Value = 168
$ws.Cells.Item(70, 5).Value = 1278
$ws.Cells.Item(70, 6).Value = 8
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 4

# Etiopia (row 139)
$ws.Cells.Item(139, 1).Value = "Etiopia"
$ws.Cells.Item(139, 2).Value = 105
$ws.Cells.Item(139, 3).Value = 9
$ws.Cells.Item(139, 4).Value = 16
$ws.Cells.Item(139, 5).Value = 86
$ws.Cells.Item(139, 6).Value = 1
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 3

# --- Zambia climbs past Cabo Verde, Polinesia Francesa, Uganda and Bahamas --
# New sort order (descending Casos totales) for rows 154-158:

# Row 154 -> Zambia
$ws.Cells.Item(154, 1).Value = "Zambia"
$ws.Cells.Item(154, 2).Value = 57
$ws.Cells.Item(154, 3).Value = 5
$ws.Cells.Item(154, 4).Value = 33
$ws.Cells.Item(154, 5).Value = 22
$ws.Cells.Item(154, 6).Value = 1
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 2

# Row 155 -> Cabo Verde
$ws.Cells.Item(155, 1).Value = "Cabo Verde"
$ws.Cells.Item(155, 2).Value = 56
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 1
$ws.Cells.Item(155, 5).Value = 54
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 1

# Row 156 -> Polinesia Francesa
$ws.Cells.Item(156, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(156, 2).Value = 55
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 53
$ws.Cells.Item(156, 6).Value = 1
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 0

# Row 157 -> Uganda
$ws.Cells.Item(157, 1).Value = "Uganda"
$ws.Cells.Item(157, 2).Value = 55
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 20
$ws.Cells.Item(157, 5).Value = 35
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0

# Row 158 -> Bahamas
$ws.Cells.Item(158, 1).Value = "Bahamas"
$ws.Cells.Item(158, 2).Value = 54
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 9
$ws.Cells.Item(158, 5).Value = 36
$ws.Cells.Item(158, 6).Value = 1
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 9

# --- Fiyi / Islas Virgenes de los Estados Unidos swap (tied data) ----------
# Row 183 -> Fiyi
$ws.Cells.Item(183, 1).Value = "Fiyi"
$ws.Cells.Item(183, 2).Value = 17
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 17
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

# Row 184 -> Islas Virgenes de los Estados Unidos
$ws.Cells.Item(184, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(184, 2).Value = 17
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 17
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

# --- Montserrat overtakes Islas Turcas y Caicos -----------------------------
# Row 195 -> Montserrat
$ws.Cells.Item(195, 1).Value = "Montserrat"
$ws.Cells.Item(195, 2).Value = 11
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 1
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = 1
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

# Row 196 -> Islas Turcas y Caicos
$ws.Cells.Item(196, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(196, 2).Value = 11
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 10
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1
